$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.447.72'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.850.48'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6302'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07700'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2934'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07746'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('D12').Value = '1.888.70'
$ws.Range('E12').Value = '  +2.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.035'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.20%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.00001074'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.89%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6795'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.73'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.84%  '
$ws.Range('D17').Value = '2.169.80'
$ws.Range('E17').Value = '  +2.86%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.198'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').Value = '29.493.58'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '228.82'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E21').Value = '  +0.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.461'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.50'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.48%  '
$ws.Range('E26').Value = '  -0.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.408'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.70'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.346'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.466'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('E31').Value = '  +0.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.129'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.037'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.849'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.10%  '
$ws.Range('E35').Value = '  +0.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7041'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.584'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01792'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('D40').Value = '1.220.43'
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.562'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9078'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.65%  '
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.85'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.36'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000121'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.137'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E48').Value = '  +0.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.004'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.686'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1147'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.52%  '
